# Applies the "se cambia los graficos" restructuring of Sheet1:
#  - header row: columns shift (Nombre column dropped, a new
#    "CARPETA DE GESTION ELECTRO" column and a trailing "Pregunta"
#    column are introduced), and
#  - every data row (2-8) has its Sucursal/checklist/comment values
#    shifted left into the new column layout, with a per-row
#    concatenated checklist-question string written into column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 ----
$r1C = $ws.Cells.Item(1,3).Value2
$r1E = $ws.Cells.Item(1,5).Value2
$r1F = $ws.Cells.Item(1,6).Value2
$r1G = $ws.Cells.Item(1,7).Value2
$r1H = $ws.Cells.Item(1,8).Value2
$r1I = $ws.Cells.Item(1,9).Value2
$ws.Cells.Item(1,2).Value = $r1C
$ws.Cells.Item(1,3).Value = "CARPETA DE GESTIÓN ELECTRO"
$ws.Cells.Item(1,4).Value = $r1E
$ws.Cells.Item(1,5).Value = $r1F
$ws.Cells.Item(1,6).Value = $r1G
$ws.Cells.Item(1,7).Value = $r1H
$ws.Cells.Item(1,8).Value = $r1I
$ws.Cells.Item(1,9).Value = "Pregunta"

# ---- Row 2 ----
$r2C = $ws.Cells.Item(2,3).Value2
$r2E = $ws.Cells.Item(2,5).Value2
$r2F = $ws.Cells.Item(2,6).Value2
$r2H = $ws.Cells.Item(2,8).Value2
$r2I = $ws.Cells.Item(2,9).Value2
$ws.Cells.Item(2,2).Value = $r2C
$ws.Cells.Item(2,3).Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene cierre y devoluciones realizadas mes anterior?;¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;Comunicación: firma de procesos claves;Tiene el Gerente realizado el check list del día?;"
$ws.Cells.Item(2,4).Value = $r2E
$ws.Cells.Item(2,5).Value = $r2F
$ws.Cells.Item(2,7).Value = $r2H
$ws.Cells.Item(2,8).Value = $r2I
$ws.Cells.Item(2,6).ClearContents()
$ws.Cells.Item(2,9).ClearContents()

# ---- Row 3 ----
$r3C = $ws.Cells.Item(3,3).Value2
$r3E = $ws.Cells.Item(3,5).Value2
$r3F = $ws.Cells.Item(3,6).Value2
$r3G = $ws.Cells.Item(3,7).Value2
$r3H = $ws.Cells.Item(3,8).Value2
$r3I = $ws.Cells.Item(3,9).Value2
$ws.Cells.Item(3,2).Value = $r3C
$ws.Cells.Item(3,3).Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;¿Tiene cierre y devoluciones realizadas mes anterior?;Comunicación: firma de procesos claves;"
$ws.Cells.Item(3,4).Value = $r3E
$ws.Cells.Item(3,5).Value = $r3F
$ws.Cells.Item(3,6).Value = $r3G
$ws.Cells.Item(3,7).Value = $r3H
$ws.Cells.Item(3,8).Value = $r3I
$ws.Cells.Item(3,9).ClearContents()

# ---- Row 4 ----
$r4C = $ws.Cells.Item(4,3).Value2
$r4E = $ws.Cells.Item(4,5).Value2
$r4F = $ws.Cells.Item(4,6).Value2
$r4G = $ws.Cells.Item(4,7).Value2
$r4H = $ws.Cells.Item(4,8).Value2
$r4I = $ws.Cells.Item(4,9).Value2
$ws.Cells.Item(4,2).Value = $r4C
$ws.Cells.Item(4,3).Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;¿Tiene cierre y devoluciones realizadas mes anterior?;Comunicación: firma de procesos claves;Tiene el Gerente realizado el check list del día?;"
$ws.Cells.Item(4,4).Value = $r4E
$ws.Cells.Item(4,5).Value = $r4F
$ws.Cells.Item(4,6).Value = $r4G
$ws.Cells.Item(4,7).Value = $r4H
$ws.Cells.Item(4,8).Value = $r4I
$ws.Cells.Item(4,9).ClearContents()

# ---- Row 5 ----
$r5C = $ws.Cells.Item(5,3).Value2
$ws.Cells.Item(5,2).Value = $r5C
$ws.Cells.Item(5,3).Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;¿Tiene cierre y devoluciones realizadas mes anterior?;Comunicación: firma de procesos claves;Tiene el Gerente realizado el check list del día?;"
$ws.Cells.Item(5,7).Value = "Ok en líneas generales_x000D_`n"
$ws.Cells.Item(5,4).ClearContents()
$ws.Cells.Item(5,8).ClearContents()

# ---- Row 6 ----
$r6C = $ws.Cells.Item(6,3).Value2
$r6H = $ws.Cells.Item(6,8).Value2
$ws.Cells.Item(6,2).Value = $r6C
$ws.Cells.Item(6,3).Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;¿Tiene cierre y devoluciones realizadas mes anterior?;Comunicación: firma de procesos claves;"
$ws.Cells.Item(6,7).Value = $r6H
$ws.Cells.Item(6,4).ClearContents()
$ws.Cells.Item(6,8).ClearContents()

# ---- Row 7 ----
$r7C = $ws.Cells.Item(7,3).Value2
$r7H = $ws.Cells.Item(7,8).Value2
$ws.Cells.Item(7,2).Value = $r7C
$ws.Cells.Item(7,3).Value = "¿Tiene firmados los objetivos de todos los vendedores?;¿Tiene planificación de trabajo por el desvío de objetivos -mes anterior? (template);¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;¿Tiene cierre y devoluciones realizadas mes anterior?;Comunicación: firma de procesos claves;Tiene el Gerente realizado el check list del día?;"
$ws.Cells.Item(7,7).Value = $r7H
$ws.Cells.Item(7,4).ClearContents()
$ws.Cells.Item(7,8).ClearContents()

# ---- Row 8 ----
$r8C = $ws.Cells.Item(8,3).Value2
$r8H = $ws.Cells.Item(8,8).Value2
$ws.Cells.Item(8,2).Value = $r8C
$ws.Cells.Item(8,3).Value = "¿Tiene acta de reuniones de los objetivos BC - BL y GEX?;Comunicación: firma de procesos claves;Tiene el Gerente realizado el check list del día?;"
$ws.Cells.Item(8,7).Value = $r8H
$ws.Cells.Item(8,4).ClearContents()
$ws.Cells.Item(8,8).ClearContents()
